# Update "想去人数" (want-to-go count) figures for a handful of events on
# the "展览" and "全部类型" sheets, per the output regenerated at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): rows 3, 7, 8, 19, 21
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1393
$ws1.Range("F7").Value = 11824
$ws1.Range("F8").Value = 4419
$ws1.Range("F19").Value = 190
$ws1.Range("F21").Value = 11373

# Sheet "全部类型" (sheet4): rows 3, 7, 8, 20, 22 (shifted by the extra
# row present only on this aggregated sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1393
$ws4.Range("F7").Value = 11824
$ws4.Range("F8").Value = 4419
$ws4.Range("F20").Value = 190
$ws4.Range("F22").Value = 11373
